# Applies the "Updated cryptos list" data refresh to Sheet1.
# Column D values that look like plain numbers are written with a
# leading apostrophe so Excel keeps them as text (matching the
# source workbook, where Price is stored as a text string, e.g. "63.930.60").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '64.292.50'
$ws.Range('E2').Value = '  -4.60%  '

# Row 3
$ws.Range('D3').Value = '3.062.26'
$ws.Range('E3').Value = '  -5.74%  '

# Row 4
$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  +0.02%  '

# Row 5
$ws.Range('D5').Value = '''562.09'
$ws.Range('E5').Value = '  -4.08%  '

# Row 6
$ws.Range('D6').Value = '''142.02'
$ws.Range('E6').Value = '  -7.00%  '

# Row 7
$ws.Range('E7').Value = '  -0.11%  '

# Row 8
$ws.Range('D8').Value = '3.053.02'
$ws.Range('E8').Value = '  -5.75%  '

# Row 9
$ws.Range('D9').Value = '''0.491'
$ws.Range('E9').Value = '  -9.83%  '

# Row 10
$ws.Range('D10').Value = '''0.158'
$ws.Range('E10').Value = '  -8.62%  '

# Row 11
$ws.Range('D11').Value = '''6.19'
$ws.Range('E11').Value = '  -9.26%  '

# Row 12
$ws.Range('D12').Value = '''0.463'
$ws.Range('E12').Value = '  -8.63%  '

# Row 13
$ws.Range('D13').Value = '''35.44'
$ws.Range('E13').Value = '  -8.13%  '

# Row 14
$ws.Range('D14').Value = '''0.0000223'
$ws.Range('E14').Value = '  -9.01%  '

# Row 15
$ws.Range('D15').Value = '3.556.82'
$ws.Range('E15').Value = '  -5.67%  '

# Row 16
$ws.Range('D16').Value = '64.303.64'
$ws.Range('E16').Value = '  -4.74%  '

# Row 17
$ws.Range('D17').Value = '''0.111'
$ws.Range('E17').Value = '  -3.15%  '

# Row 18
$ws.Range('D18').Value = '3.068.06'
$ws.Range('E18').Value = '  -5.65%  '

# Row 19
$ws.Range('B19').Value = 'BitcoinCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D19').Value = '''488.05'
$ws.Range('E19').Value = '  -10.46%  '

# Row 20
$ws.Range('B20').Value = 'Polkadot'
$ws.Range('C20').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D20').Value = '''6.60'
$ws.Range('E20').Value = '  -9.00%  '

# Row 21
$ws.Range('D21').Value = '''13.68'
$ws.Range('E21').Value = '  -10.24%  '

# Row 22
$ws.Range('D22').Value = '''0.682'
$ws.Range('E22').Value = '  -11.06%  '

# Row 23
$ws.Range('D23').Value = '''7.20'
$ws.Range('E23').Value = '  -8.44%  '

# Row 24
$ws.Range('D24').Value = '''12.47'
$ws.Range('E24').Value = '  -7.66%  '

# Row 25
$ws.Range('D25').Value = '''78.08'
$ws.Range('E25').Value = '  -8.81%  '

# Row 26
$ws.Range('E26').Value = '  +0.02%  '

# Row 27
$ws.Range('D27').Value = '''2.76'
$ws.Range('E27').Value = '  -14.05%  '

# Row 28
$ws.Range('B28').Value = 'RenderToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D28').Value = '''7.81'
$ws.Range('E28').Value = '  -4.95%  '

# Row 29
$ws.Range('B29').Value = 'ImmutableX'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D29').Value = '''2.07'
$ws.Range('E29').Value = '  -3.14%  '

# Row 30
$ws.Range('E30').Value = '  -0.17%  '

# Row 31
$ws.Range('D31').Value = '''26.40'
$ws.Range('E31').Value = '  -10.37%  '

# Row 32
$ws.Range('D32').Value = '''2.61'
$ws.Range('E32').Value = '  -4.30%  '

# Row 33
$ws.Range('D33').Value = '''1.12'
$ws.Range('E33').Value = '  -3.29%  '

# Row 34
$ws.Range('D34').Value = '''509.67'
$ws.Range('E34').Value = '  -7.95%  '

# Row 35
$ws.Range('D35').Value = '''5.45'
$ws.Range('E35').Value = '  -5.55%  '

# Row 36
$ws.Range('D36').Value = '''52.94'
$ws.Range('E36').Value = '  -1.65%  '

# Row 37
$ws.Range('D37').Value = '''5.89'
$ws.Range('E37').Value = '  -10.85%  '

# Row 38
$ws.Range('D38').Value = '''0.0406'
$ws.Range('E38').Value = '  -8.08%  '

# Row 39
$ws.Range('D39').Value = '''0.0791'
$ws.Range('E39').Value = '  -7.01%  '

# Row 40
$ws.Range('D40').Value = '''0.120'
$ws.Range('E40').Value = '  -6.17%  '

# Row 41
$ws.Range('D41').Value = '''8.30'
$ws.Range('E41').Value = '  -10.15%  '

# Row 42
$ws.Range('B42').Value = 'dogwifhat'
$ws.Range('C42').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D42').Value = '''2.67'
$ws.Range('E42').Value = '  +2.05%  '

# Row 43
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '2.879.22'
$ws.Range('E43').Value = '  -2.10%  '

# Row 44
$ws.Range('E44').Value = '  -0.22%  '

# Row 45
$ws.Range('D45').Value = '''0.243'
$ws.Range('E45').Value = '  -7.41%  '

# Row 46
$ws.Range('D46').Value = '0.0₃0544'
$ws.Range('E46').Value = '  -6.93%  '

# Row 47
$ws.Range('D47').Value = '''2.06'
$ws.Range('E47').Value = '  -3.50%  '

# Row 48
$ws.Range('D48').Value = '''24.73'
$ws.Range('E48').Value = '  -6.29%  '

# Row 49
$ws.Range('D49').Value = '''118.38'
$ws.Range('E49').Value = '  -5.78%  '

# Row 50
$ws.Range('D50').Value = '''0.107'
$ws.Range('E50').Value = '  -5.78%  '

# Row 51
$ws.Range('D51').Value = '''2.06'
$ws.Range('E51').Value = '  -13.25%  '
